$d = $word.ActiveDocument

# --- Step 1: replace the final paragraph with the 3-paragraph block -------
# (reformatted paragraph + 2 new paragraphs) via a whole-paragraph InsertXML.
$target = $d.Paragraphs.Last
$r = $target.Range
$r2 = $d.Range($r.Start, $r.End - 1)

$fragment = @'
<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="4D4D4D"/><w:spacing w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:shd w:val="clear" w:fill="FFFFFF"/></w:rPr><w:t>因为当前sass的版本太高</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="4D4D4D"/><w:spacing w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/></w:rPr><w:t>直接用npm install sass-loader@7.3.1 --save-dev 就可以了，同名的第三库会自动更换</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:ind w:left="0" w:leftChars="0" w:firstLine="0" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="default" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Entrypoint mini-css-extract-plugin = * </w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei" w:cs="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei"/><w:i w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t>W</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Microsoft YaHei" w:hAnsi="Microsoft YaHei" w:eastAsia="Microsoft YaHei"/><w:i w:val="0"/><w:caps w:val="0"/><w:color w:val="222226"/><w:spacing w:val="0"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:shd w:val="clear" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t>hen using miniCssExtractPlugin, it appears the above yellow text. To remove it, uninstall miniCssExtractPlugin, only use style-loader.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2.InsertXML($xml)

# --- Step 2: restore explicit "off" toggles (w:i val=0 / w:caps val=0) ----
# A whole-paragraph InsertXML silently drops explicit-false boolean toggle
# properties: force them back on by bouncing the COM property true->false,
# which makes Word re-emit the explicit <w:i w:val="0"/> / <w:caps w:val="0"/>.
$n = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($n - 2)
$p2 = $d.Paragraphs.Item($n - 1)
$p3 = $d.Paragraphs.Item($n)

# Paragraphs 1 and 2 (of the 3): every run + the paragraph mark need both
# Italic=0 and AllCaps=0 explicitly.
foreach ($p in @($p1, $p2)) {
  $pr = $p.Range
  $pr.Font.Italic = 1
  $pr.Font.Italic = 0
  $pr.Font.AllCaps = 1
  $pr.Font.AllCaps = 0
}

# Paragraph 3: Italic=0 applies to everything (runs + mark), but AllCaps=0
# must skip the leading "W" run (it has no w:caps at all in the target).
$pr3 = $p3.Range
$pr3.Font.Italic = 1
$pr3.Font.Italic = 0

$wRun = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
Write-Host "wRun text:" $wRun.Text
$restRange = $d.Range($p3.Range.Start + 1, $p3.Range.End)
$restRange.Font.AllCaps = 1
$restRange.Font.AllCaps = 0

Write-Host "Paragraphs after fixups:" $d.Paragraphs.Count
